# Refresh cryptocurrency snapshot data (Price, Volume(1h), Hora columns)
# per the "Updated symbol list" GitHub Actions commit.
# Values must remain plain text (they were authored as inline strings,
# e.g. "277.80", "1.10%", "14"), so we force NumberFormat to Text before
# writing and reset the style back to Normal afterwards so no stray
# number formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "277.80"
Set-TextValue "E2" "1.10%"
Set-TextValue "G2" "14"
Set-TextValue "D3" "27.19"
Set-TextValue "E3" "-0.72%"
Set-TextValue "G3" "14"
Set-TextValue "D4" "4.858"
Set-TextValue "E4" "1.37%"
Set-TextValue "G4" "14"
Set-TextValue "D5" "0.06394"
Set-TextValue "E5" "1.57%"
Set-TextValue "G5" "14"
Set-TextValue "D6" "6.998"
Set-TextValue "E6" "1.13%"
Set-TextValue "G6" "14"
Set-TextValue "E7" "-6.83%"
Set-TextValue "G7" "14"
Set-TextValue "D8" "0.8789"
Set-TextValue "E8" "0.78%"
Set-TextValue "G8" "14"
Set-TextValue "D9" "0.1522"
Set-TextValue "E9" "0.70%"
Set-TextValue "G9" "14"
Set-TextValue "D10" "0.05177"
Set-TextValue "E10" "3.45%"
Set-TextValue "G10" "14"
Set-TextValue "D11" "0.07509"
Set-TextValue "E11" "1.21%"
Set-TextValue "G11" "14"
Set-TextValue "E12" "1.26%"
Set-TextValue "G12" "14"
Set-TextValue "D13" "0.08965"
Set-TextValue "G13" "14"
Set-TextValue "D14" "0.001570"
Set-TextValue "E14" "0.64%"
Set-TextValue "G14" "14"
Set-TextValue "D15" "0.0006386"
Set-TextValue "E15" "0.72%"
Set-TextValue "G15" "14"
Set-TextValue "D16" "0.006076"
Set-TextValue "E16" "2.73%"
Set-TextValue "G16" "14"
Set-TextValue "D17" "3.476"
Set-TextValue "E17" "0.75%"
Set-TextValue "G17" "14"
Set-TextValue "D18" "3.300"
Set-TextValue "E18" "-0.13%"
Set-TextValue "G18" "14"
Set-TextValue "D19" "2.248"
Set-TextValue "E19" "-1.59%"
Set-TextValue "G19" "14"
Set-TextValue "D20" "0.3141"
Set-TextValue "E20" "-0.23%"
Set-TextValue "G20" "14"
Set-TextValue "D21" "0.1324"
Set-TextValue "E21" "0.45%"
Set-TextValue "G21" "14"
Set-TextValue "D22" "3.903"
Set-TextValue "E22" "0.14%"
Set-TextValue "G22" "14"
Set-TextValue "D23" "0.04412"
Set-TextValue "E23" "0.72%"
Set-TextValue "G23" "14"
Set-TextValue "D24" "0.1505"
Set-TextValue "E24" "9.03%"
Set-TextValue "G24" "14"
Set-TextValue "D25" "0.001175"
Set-TextValue "E25" "0.60%"
Set-TextValue "G25" "14"
Set-TextValue "D26" "0.003897"
Set-TextValue "E26" "1.95%"
Set-TextValue "G26" "14"
Set-TextValue "G27" "14"
Set-TextValue "D28" "0.0001180"
Set-TextValue "E28" "-1.58%"
Set-TextValue "G28" "14"
Set-TextValue "E29" "1.69%"
Set-TextValue "G29" "14"
Set-TextValue "G30" "14"
Set-TextValue "G31" "14"
Set-TextValue "G32" "14"
Set-TextValue "G33" "14"
Set-TextValue "G34" "14"
Set-TextValue "G35" "14"
Set-TextValue "G36" "14"
Set-TextValue "G37" "14"
Set-TextValue "G38" "14"
Set-TextValue "G39" "14"
Set-TextValue "D40" "0.04075"
Set-TextValue "E40" "-0.41%"
Set-TextValue "G40" "14"
Set-TextValue "D41" "0.006889"
Set-TextValue "E41" "-2.13%"
Set-TextValue "G41" "14"
Set-TextValue "E42" "0.26%"
Set-TextValue "G42" "14"
Set-TextValue "E43" "-6.40%"
Set-TextValue "G43" "14"
Set-TextValue "D44" "0.01123"
Set-TextValue "E44" "0.19%"
Set-TextValue "G44" "14"
Set-TextValue "D45" "0.00005355"
Set-TextValue "E45" "3.26%"
Set-TextValue "G45" "14"
Set-TextValue "E46" "5.01%"
Set-TextValue "G46" "14"
Set-TextValue "D47" "0.01851"
Set-TextValue "E47" "-19.45%"
Set-TextValue "G47" "14"
Set-TextValue "G48" "14"
Set-TextValue "G49" "14"
Set-TextValue "G50" "14"
Set-TextValue "G51" "14"
